$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.588.20"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.585.46"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.810.46"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "1.578.39"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "27.570.02"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "1.375.50"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.975"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.829"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("D47").Value = "1.720.96"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "0.0₇0992"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0496"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
